# R Code for Sensor and Master tables
# R Code for Sensor and wl with worldview 3 actualized.
#
# Applies the "BandsSensors" sensor-table update to Folha3 (sheet2):
#  - inserts 5 new SWIR band columns (SWIR3..SWIR7, min/max) after SWIR2,
#    shifting the old TIRS1/TIRS2 min/max columns further right
#  - fills in the newly created columns for the Worldview row
#  - splits "Worldview-2 and 3" into two separate rows: "Worldview-2"
#    (row 12, existing row updated) and "Worldview-3" (row 13, new row)
#  - adds a couple of previously-missing VRE2 values on row 11/12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha3")

# --- Row 1 (header): insert 10 new headers (SWIR3_min .. SWIR7_max) right
#     after SWIR2_max (AD1), then re-append the TIRS headers that used to
#     sit at AE1:AH1 into their new slots AO1:AR1.
$ws.Range("AE1").Value = "SWIR3_min"
$ws.Range("AF1").Value = "SWIR3_max"
$ws.Range("AG1").Value = "SWIR4_min"
$ws.Range("AH1").Value = "SWIR4_max"
$ws.Range("AI1").Value = "SWIR5_min"
$ws.Range("AJ1").Value = "SWIR5_max"
$ws.Range("AK1").Value = "SWIR6_min"
$ws.Range("AL1").Value = "SWIR6_max"
$ws.Range("AM1").Value = "SWIR7_min"
$ws.Range("AN1").Value = "SWIR7_max"
$ws.Range("AO1").Value = "TIRS1_min"
$ws.Range("AP1").Value = "TIRS1_max"
$ws.Range("AQ1").Value = "TIRS2_min"
$ws.Range("AR1").Value = "TIRS2_max"

# --- Row 2 (Sentinel-2): move the SWIR2/TIRS1/TIRS2 values into their new
#     columns, then clear the now-unused old cells.
$ws.Range("AG2").Value = 2.0720000000000001
$ws.Range("AH2").Value = 2.3119999999999998
$ws.Range("AO2").Value = 10.6
$ws.Range("AP2").Value = 11.19
$ws.Range("AQ2").Value = 11.5
$ws.Range("AR2").Value = 12.51
$ws.Range("AC2").Value = ""
$ws.Range("AD2").Value = ""
$ws.Range("AE2").Value = ""
$ws.Range("AF2").Value = ""

# --- Row 3 (Landsat-8): same column shift.
$ws.Range("AG3").Value = 2.1070000000000002
$ws.Range("AH3").Value = 2.294
$ws.Range("AO3").Value = 10.6
$ws.Range("AP3").Value = 11.19
$ws.Range("AQ3").Value = 11
$ws.Range("AR3").Value = 12.005000000000001
$ws.Range("AC3").Value = ""
$ws.Range("AD3").Value = ""
$ws.Range("AE3").Value = ""
$ws.Range("AF3").Value = ""

# --- Row 11 (Kompsat 2, 3 and 3A): fill in the VRE2 min/max that were
#     previously blank.
$ws.Range("O11").Value = 0.70499999999999996
$ws.Range("P11").Value = 0.745

# --- Row 12: this used to be the combined "Worldview-2 and 3" row; it now
#     becomes the "Worldview-2"-only row, with corrected Yellow values and
#     newly populated NIR2/WV/Cirrus/SWIR1/SWIR2..SWIR7 columns.
$ws.Range("A12").Value = "Worldview-2"
$ws.Range("I12").Value = 0.58499999999999996
$ws.Range("J12").Value = 0.625
$ws.Range("S12").Value = 0.77
$ws.Range("T12").Value = 0.89500000000000002
$ws.Range("U12").Value = 0.86
$ws.Range("V12").Value = 1.04
$ws.Range("Y12").Value = 1.1950000000000001
$ws.Range("Z12").Value = 1.2250000000000001
$ws.Range("AA12").Value = 1.55
$ws.Range("AB12").Value = 1.59
$ws.Range("AC12").Value = 1.64
$ws.Range("AD12").Value = 1.68
$ws.Range("AE12").Value = 1.71
$ws.Range("AF12").Value = 1.75
$ws.Range("AG12").Value = 2.145
$ws.Range("AH12").Value = 2.1850000000000001
$ws.Range("AI12").Value = 2.1850000000000001
$ws.Range("AJ12").Value = 2.2250000000000001
$ws.Range("AK12").Value = 2.2349999999999999
$ws.Range("AL12").Value = 2.2850000000000001
$ws.Range("AM12").Value = 2.2949999999999999
$ws.Range("AN12").Value = 2.3650000000000002

# --- Row 13 (new): "Worldview-3" row.
$ws.Range("A13").Value = "Worldview-3"
$ws.Range("B13").Value = "1.24,3.7"
$ws.Range("C13").Value = 0.4
$ws.Range("D13").Value = 0.45
$ws.Range("E13").Value = 0.45
$ws.Range("F13").Value = 0.51
$ws.Range("G13").Value = 0.51
$ws.Range("H13").Value = 0.57999999999999996
$ws.Range("I13").Value = 0.57999999999999996
$ws.Range("K13").Value = 0.63
$ws.Range("L13").Value = 0.69
$ws.Range("O13").Value = 0.70499999999999996
$ws.Range("P13").Value = 0.745
$ws.Range("S13").Value = 0.77
$ws.Range("T13").Value = 0.89500000000000002
$ws.Range("U13").Value = 0.86
$ws.Range("V13").Value = 1.04

# --- Selection: mirror the saved workbook's final selection over the full
#     used range of the table.
$ws.Range("A1:AR13").Select()
